$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.842.75'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.639.56'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.383'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.78%  '
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '3.115.78'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Value = '63.700.16'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('E16').Value = '  +1.47%  '
$ws.Range('D17').Value = '2.653.45'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.79%  '
$ws.Range('E19').Value = '  +3.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -1.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.65%  '
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('E27').Value = '  +6.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '565.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '169.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.31%  '
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '165.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0568'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0246'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0956'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
